# Apply weekly Fruta/Hortaliza data update to the "Pomelo" sheet.
# The diff shows the date/volume/price data for rows 3, 4, 6, 8 and 9
# being rotated among themselves (two cycles: 3->6->9->3 and 4->8->4),
# i.e. each row ends up with another row's D/M/N/O/P/Q/S values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44309
$ws.Range("M3").Value = 300
$ws.Range("N3").Value = 7000
$ws.Range("O3").Value = 7000
$ws.Range("P3").Value = 7000
$ws.Range("S3").Value = 500

# Row 4
$ws.Range("D4").Value = 44176
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = 7000
$ws.Range("O4").Value = 7000
$ws.Range("P4").Value = 7000
$ws.Range("Q4").Value = "$/caja 14 kilos empedrada"
$ws.Range("S4").Value = 500

# Row 6
$ws.Range("D6").Value = 44162
$ws.Range("M6").Value = 120

# Row 8
$ws.Range("D8").Value = 44397
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 11000
$ws.Range("O8").Value = 11000
$ws.Range("P8").Value = 11000
$ws.Range("Q8").Value = "$/caja 14 kilos"
$ws.Range("S8").Value = 786

# Row 9
$ws.Range("D9").Value = 44208
$ws.Range("M9").Value = 210
$ws.Range("N9").Value = 10000
$ws.Range("O9").Value = 10000
$ws.Range("P9").Value = 10000
$ws.Range("S9").Value = 714
